# Update cryptos list values (price / 1h volume, plus a few reordered rows)
# as produced by the scheduled GitHub Actions scraper run.
#
# Note: several "Price" values in column D look like plain decimal numbers
# (e.g. "0.07502", "42.48"). Excel's COM layer auto-converts such strings to
# numeric cell values, which does not match the source data (plain text).
# To force them to stay text we prefix with an apostrophe (Excel's classic
# "treat as text" marker) and then reset the cell style to "Normal" so no
# stray number-format / quote-prefix style remains attached to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "28.638.20"

$ws.Cells.Item(3, 4).Value = "1.802.46"
$ws.Cells.Item(3, 5).Value = "  -0.76%  "

$ws.Cells.Item(4, 5).Value = "  +0.14%  "

$ws.Cells.Item(5, 4).Value = "'316.66"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.55%  "

$ws.Cells.Item(6, 5).Value = "  +0.19%  "

$ws.Cells.Item(7, 4).Value = "'0.5311"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -7.82%  "

$ws.Cells.Item(8, 4).Value = "'0.3766"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -2.61%  "

$ws.Cells.Item(9, 2).Value = "Dogecoin"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(9, 4).Value = "'0.07502"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -1.60%  "

$ws.Cells.Item(10, 2).Value = "OKB"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(10, 4).Value = "'42.48"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -1.55%  "

$ws.Cells.Item(11, 5).Value = "  -2.28%  "

$ws.Cells.Item(12, 5).Value = "  +0.19%  "

$ws.Cells.Item(13, 4).Value = "'20.68"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -2.63%  "

$ws.Cells.Item(14, 4).Value = "'6.152"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -1.72%  "

$ws.Cells.Item(15, 4).Value = "'7.348"
$ws.Cells.Item(15, 4).Style = "Normal"

$ws.Cells.Item(16, 4).Value = "1.799.14"
$ws.Cells.Item(16, 5).Value = "  -0.78%  "

$ws.Cells.Item(17, 4).Value = "'90.22"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -2.10%  "

$ws.Cells.Item(18, 4).Value = "'0.00001064"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -1.23%  "

$ws.Cells.Item(19, 4).Value = "'0.06465"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.14%  "

$ws.Cells.Item(20, 5).Value = "  +0.20%  "

$ws.Cells.Item(21, 4).Value = "'17.24"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.48%  "

$ws.Cells.Item(22, 4).Value = "'5.903"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -1.45%  "

$ws.Cells.Item(23, 4).Value = "28.659.25"
$ws.Cells.Item(23, 5).Value = "  +0.98%  "

$ws.Cells.Item(24, 4).Value = "'11.12"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -1.69%  "

$ws.Cells.Item(25, 4).Value = "'2.101"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.05%  "

$ws.Cells.Item(26, 4).Value = "'159.21"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +0.90%  "

$ws.Cells.Item(27, 4).Value = "'20.45"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -1.88%  "

$ws.Cells.Item(28, 4).Value = "2.008.86"
$ws.Cells.Item(28, 5).Value = "  -0.78%  "

$ws.Cells.Item(29, 4).Value = "'2.346"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -2.62%  "

$ws.Cells.Item(30, 4).Value = "'122.76"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -0.66%  "

$ws.Cells.Item(31, 5).Value = "  -5.93%  "

$ws.Cells.Item(32, 4).Value = "'0.1057"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -0.69%  "

$ws.Cells.Item(33, 5).Value = "  -2.28%  "

$ws.Cells.Item(34, 4).Value = "'3.686"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +1.54%  "

$ws.Cells.Item(35, 4).Value = "'0.2235"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +3.02%  "

$ws.Cells.Item(36, 4).Value = "'0.06394"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +5.09%  "

$ws.Cells.Item(37, 4).Value = "'0.02308"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.32%  "

$ws.Cells.Item(38, 4).Value = "'8.770"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -1.40%  "

$ws.Cells.Item(39, 4).Value = "'5.042"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +0.10%  "

$ws.Cells.Item(40, 2).Value = "Aptos"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(40, 4).Value = "'11.24"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -3.74%  "

$ws.Cells.Item(41, 2).Value = "TrustWalletToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(41, 4).Value = "'1.202"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +4.01%  "

$ws.Cells.Item(42, 4).Value = "'0.6200"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -3.36%  "

$ws.Cells.Item(43, 2).Value = "WEMIXTOKEN"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(43, 4).Value = "'1.416"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +2.65%  "

$ws.Cells.Item(44, 2).Value = "Frax"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(44, 4).Value = "'1.002"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.19%  "

$ws.Cells.Item(45, 5).Value = "  -1.90%  "

$ws.Cells.Item(46, 4).Value = "'3.691"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -0.36%  "

$ws.Cells.Item(47, 4).Value = "'0.5843"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -2.34%  "

$ws.Cells.Item(48, 4).Value = "'125.89"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +3.11%  "

$ws.Cells.Item(49, 4).Value = "'1.939"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -0.14%  "

$ws.Cells.Item(50, 4).Value = "'1.151"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.36%  "

$ws.Cells.Item(51, 4).Value = "'0.06891"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.68%  "
